# Append " (Changed main)" to the end of the first paragraph
# ("This is a Microsoft word document."), split across three new
# runs: " (", "Changed main", ")" -- matching the target OOXML diff.
#
# A plain Range.InsertAfter() would work for the visible text, but this
# engine (like Word itself) silently coalesces adjacent runs that end up
# with identical formatting once the document is saved, which would
# collapse the four runs back into one. To keep the runs distinct (and
# free of any leftover direct-formatting markup) we rebuild the whole
# paragraph's content via Range.InsertXML() using a minimal flat-OPC
# WordprocessingML package: InsertXML() *replaces* the exact range it is
# called on, so targeting the full paragraph content (excluding the
# trailing paragraph mark) swaps it in-place for the new four-run
# version without disturbing any other paragraph, paragraph id, or
# section property.

$d = $word.ActiveDocument

# Locate the paragraph to edit by its current text, rather than a
# hard-coded character offset, so the script is resilient to unrelated
# edits elsewhere in the document.
$found = $d.Content
$found.Find.Execute("This is a Microsoft word document.", $false, $false,
                     $false, $false, $false, $true, 1, $false, "", 0)

$para = $found.Paragraphs(1)
$paraStart = $para.Range.Start
$paraEnd = $para.Range.End - 1   # exclude the trailing paragraph mark
$target = $d.Range($paraStart, $paraEnd)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>This is a Microsoft word document.</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:t>Changed main</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xml)
